$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row labels (row 1) to reflect the new "2-Room/3-Room always offered"
# columns instead of the generic "Type 1 / Type 2" + units columns.
$ws.Range("C1").Value = "Always 2-Room"
$ws.Range("D1").Value = "Number of 2-Room"
$ws.Range("F1").Value = "Always 3-Room"
$ws.Range("G1").Value = "Number of 3-Room"

# Update the active selection to match the latest saved view.
$ws.Range("E7").Select()
